$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.442694664001465
$ws.Range("B1").Value = 1.402629613876343
$ws.Range("C1").Value = 3.742656469345093
$ws.Range("D1").Value = 5.707573890686035
$ws.Range("E1").Value = 1.635126233100891
